$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1733.1666
$ws.Range("I43").Value = 1349.5
$ws.Range("J43").Value = 1925
$ws.Range("K43").Value = 1349.5
$ws.Range("L43").Value = 1925
$ws.Range("M43").Value = -1280.5
$ws.Range("N43").Value = -2063
$ws.Range("H111").Value = 73492.92999999999
$ws.Range("I111").Value = 2240
$ws.Range("J111").Value = 251625.25
$ws.Range("K111").Value = 6720
$ws.Range("L111").Value = 754875.75
$ws.Range("M111").Value = -3653
$ws.Range("N111").Value = -761009.75
$ws.Range("H129").Value = 965.19446
$ws.Range("I129").Value = 881.8333
$ws.Range("J129").Value = 972.7727
$ws.Range("K129").Value = 2645.4999
$ws.Range("L129").Value = 2918.3181
$ws.Range("M129").Value = 2354.5001
$ws.Range("N129").Value = -12918.3181
$ws.Range("H132").Value = 1821.8286
$ws.Range("I132").Value = 1528.0869
$ws.Range("J132").Value = 2384.8333
$ws.Range("K132").Value = 4584.2607
$ws.Range("L132").Value = 7154.499899999999
$ws.Range("M132").Value = -2054.2607
$ws.Range("N132").Value = -12214.4999
$ws.Range("H137").Value = 2258.3125
$ws.Range("I137").Value = 1811.25
$ws.Range("J137").Value = 2705.375
$ws.Range("K137").Value = 5433.75
$ws.Range("L137").Value = 8116.125
$ws.Range("M137").Value = -2883.75
$ws.Range("N137").Value = -13216.125
$ws.Range("H138").Value = 2117.753
$ws.Range("I138").Value = 820.9583
$ws.Range("J138").Value = 3800.081
$ws.Range("K138").Value = 2462.8749
$ws.Range("L138").Value = 11400.243
$ws.Range("M138").Value = 2677.1251
$ws.Range("N138").Value = -21680.243
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4671.44
$ws.Range("I32").Value = 3355.3684
$ws.Range("J32").Value = 8839
$ws.Range("K32").Value = 3355.3684
$ws.Range("L32").Value = 8839
$ws.Range("M32").Value = -3068.3684
$ws.Range("N32").Value = -9413
$ws.Range("H44").Value = 17993.334
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 17993.334
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 17993.334
$ws.Range("N44").Value = -18969.334
$ws.Range("H55").Value = 18998.666
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 18998.666
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 18998.666
$ws.Range("N55").Value = -19628.666
$ws.Range("H61").Value = 1704.9656
$ws.Range("I61").Value = 1877.25
$ws.Range("J61").Value = 1492.9231
$ws.Range("K61").Value = 1877.25
$ws.Range("L61").Value = 1492.9231
$ws.Range("M61").Value = -1665.25
$ws.Range("N61").Value = -1916.9231
$ws.Range("H74").Value = 783.8372000000001
$ws.Range("I74").Value = 563.8570999999999
$ws.Range("J74").Value = 1746.25
$ws.Range("K74").Value = 563.8570999999999
$ws.Range("L74").Value = 1746.25
$ws.Range("M74").Value = 310.1429000000001
$ws.Range("N74").Value = -3494.25
$ws.Range("H77").Value = 783.8372000000001
$ws.Range("I77").Value = 563.8570999999999
$ws.Range("J77").Value = 1746.25
$ws.Range("K77").Value = 2819.2855
$ws.Range("L77").Value = 8731.25
$ws.Range("M77").Value = 1548.7145
$ws.Range("N77").Value = -17467.25
$ws.Range("H88").Value = 2035.1666
$ws.Range("I88").Value = 1916.8572
$ws.Range("J88").Value = 2200.8
$ws.Range("K88").Value = 1916.8572
$ws.Range("L88").Value = 2200.8
$ws.Range("M88").Value = -1510.8572
$ws.Range("N88").Value = -3012.8
$ws.Range("H91").Value = 2035.1666
$ws.Range("I91").Value = 1916.8572
$ws.Range("J91").Value = 2200.8
$ws.Range("K91").Value = 1916.8572
$ws.Range("L91").Value = 2200.8
$ws.Range("M91").Value = -512.8571999999999
$ws.Range("N91").Value = -5008.8
$ws.Range("H122").Value = 2850199.5
$ws.Range("I122").Value = 3205849.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9617548.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -9615098.5
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2328614.8
$ws.Range("I132").Value = 1977.1333
$ws.Range("J132").Value = 7697778
$ws.Range("K132").Value = 5931.3999
$ws.Range("L132").Value = 23093334
$ws.Range("M132").Value = -3401.3999
$ws.Range("N132").Value = -23098394
$ws.Range("H134").Value = 40000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 40000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140
$ws.Range("H136").Value = 1704.9656
$ws.Range("I136").Value = 1877.25
$ws.Range("J136").Value = 1492.9231
$ws.Range("K136").Value = 5631.75
$ws.Range("L136").Value = 4478.7693
$ws.Range("M136").Value = -3081.75
$ws.Range("N136").Value = -9578.7693
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 14100
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 14100
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 14100
$ws.Range("N46").Value = -14696
$ws.Range("M46").ClearContents()
$ws.Range("H86").Value = 1727.2778
$ws.Range("I86").Value = 1682.6666
$ws.Range("J86").Value = 1816.5
$ws.Range("K86").Value = 1682.6666
$ws.Range("L86").Value = 1816.5
$ws.Range("M86").Value = -559.6666
$ws.Range("N86").Value = -4062.5
$ws.Range("H89").Value = 1727.2778
$ws.Range("I89").Value = 1682.6666
$ws.Range("J89").Value = 1816.5
$ws.Range("K89").Value = 8413.333000000001
$ws.Range("L89").Value = 9082.5
$ws.Range("M89").Value = -2797.333000000001
$ws.Range("N89").Value = -20314.5
$ws.Range("H134").Value = 1839.5193
$ws.Range("I134").Value = 1497.1666
$ws.Range("J134").Value = 2306.3635
$ws.Range("K134").Value = 4491.4998
$ws.Range("L134").Value = 6919.0905
$ws.Range("M134").Value = -1956.4998
$ws.Range("N134").Value = -11989.0905
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6450.315
$ws.Range("I31").Value = 1620.125
$ws.Range("J31").Value = 20250.857
$ws.Range("K31").Value = 1620.125
$ws.Range("L31").Value = 20250.857
$ws.Range("M31").Value = -1325.125
$ws.Range("N31").Value = -20840.857
$ws.Range("H34").Value = 6450.315
$ws.Range("I34").Value = 1620.125
$ws.Range("J34").Value = 20250.857
$ws.Range("K34").Value = 1620.125
$ws.Range("L34").Value = 20250.857
$ws.Range("M34").Value = -1418.125
$ws.Range("N34").Value = -20654.857
$ws.Range("H62").Value = 8024.25
$ws.Range("I62").Value = 11992.5
$ws.Range("J62").Value = 4056
$ws.Range("K62").Value = 11992.5
$ws.Range("L62").Value = 4056
$ws.Range("M62").Value = -11368.5
$ws.Range("N62").Value = -5304
$ws.Range("H65").Value = 8024.25
$ws.Range("I65").Value = 11992.5
$ws.Range("J65").Value = 4056
$ws.Range("K65").Value = 59962.5
$ws.Range("L65").Value = 20280
$ws.Range("M65").Value = -56842.5
$ws.Range("N65").Value = -26520
$ws.Range("H132").Value = 2434.3408
$ws.Range("I132").Value = 2036.3704
$ws.Range("J132").Value = 3066.4119
$ws.Range("K132").Value = 6109.1112
$ws.Range("L132").Value = 9199.235700000001
$ws.Range("M132").Value = -3579.1112
$ws.Range("N132").Value = -14259.2357
$ws.Range("H134").Value = 2525.6904
$ws.Range("I134").Value = 3222.3462
$ws.Range("J134").Value = 1393.625
$ws.Range("K134").Value = 9667.0386
$ws.Range("L134").Value = 4180.875
$ws.Range("M134").Value = -7132.0386
$ws.Range("N134").Value = -9250.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 585
$ws.Range("I29").Value = 176
$ws.Range("J29").Value = 877.1429000000001
$ws.Range("K29").Value = 528
$ws.Range("L29").Value = 2631.4287
$ws.Range("M29").Value = -251
$ws.Range("N29").Value = -3185.4287
$ws.Range("H42").Value = 1500
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 1500
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 4500
$ws.Range("N42").Value = -5568
$ws.Range("H117").Value = 17552616
$ws.Range("I117").Value = 20325.8
$ws.Range("J117").Value = 23814148
$ws.Range("K117").Value = 60977.39999999999
$ws.Range("L117").Value = 71442444
$ws.Range("M117").Value = -57535.39999999999
$ws.Range("N117").Value = -71449328
$ws.Range("H121").Value = 928.9535
$ws.Range("I121").Value = 790
$ws.Range("J121").Value = 935.7317
$ws.Range("K121").Value = 2370
$ws.Range("L121").Value = 2807.1951
$ws.Range("M121").Value = -1060
$ws.Range("N121").Value = -5427.1951
$ws.Range("H124").Value = 6455.533
$ws.Range("I124").Value = 2000
$ws.Range("J124").Value = 6773.7856
$ws.Range("K124").Value = 6000
$ws.Range("L124").Value = 20321.3568
$ws.Range("M124").Value = -1090
$ws.Range("N124").Value = -30141.3568
$ws.Range("H129").Value = 33334948
$ws.Range("I129").Value = 66667492
$ws.Range("J129").Value = 2406.6
$ws.Range("K129").Value = 200002476
$ws.Range("L129").Value = 7219.799999999999
$ws.Range("M129").Value = -199997476
$ws.Range("N129").Value = -17219.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 11035.2
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 11035.2
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 11035.2
$ws.Range("N57").Value = -12675.2
$ws.Range("H113").Value = 34483684
$ws.Range("I113").Value = 62500652
$ws.Range("J113").Value = 1264.6154
$ws.Range("K113").Value = 62500652
$ws.Range("L113").Value = 1264.6154
$ws.Range("M113").Value = -62498482
$ws.Range("N113").Value = -5604.6154
$ws.Range("H126").Value = 4543.773
$ws.Range("I126").Value = 8433.467000000001
$ws.Range("J126").Value = 2531.862
$ws.Range("K126").Value = 25300.401
$ws.Range("L126").Value = 7595.586
$ws.Range("M126").Value = -22830.401
$ws.Range("N126").Value = -12535.586
$ws.Range("H132").Value = 2026
$ws.Range("I132").Value = 1718.2106
$ws.Range("J132").Value = 2318.4
$ws.Range("K132").Value = 5154.6318
$ws.Range("L132").Value = 6955.200000000001
$ws.Range("M132").Value = -2624.6318
$ws.Range("N132").Value = -12015.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1499.0869
$ws.Range("I61").Value = 1533.5714
$ws.Range("J61").Value = 1445.4445
$ws.Range("K61").Value = 1533.5714
$ws.Range("L61").Value = 1445.4445
$ws.Range("M61").Value = -1331.5714
$ws.Range("N61").Value = -1849.4445
$ws.Range("H113").Value = 1499.0869
$ws.Range("I113").Value = 1533.5714
$ws.Range("J113").Value = 1445.4445
$ws.Range("K113").Value = 1533.5714
$ws.Range("L113").Value = 1445.4445
$ws.Range("M113").Value = 636.4286
$ws.Range("N113").Value = -5785.4445
$ws.Range("H132").Value = 8907823
$ws.Range("I132").Value = 10688592
$ws.Range("J132").Value = 3979.2
$ws.Range("K132").Value = 32065776
$ws.Range("L132").Value = 11937.6
$ws.Range("M132").Value = -32063246
$ws.Range("N132").Value = -16997.6
$ws.Range("H136").Value = 3752.0557
$ws.Range("I136").Value = 1565.9656
$ws.Range("J136").Value = 12808.714
$ws.Range("K136").Value = 4697.8968
$ws.Range("L136").Value = 38426.142
$ws.Range("M136").Value = -2147.8968
$ws.Range("N136").Value = -43526.142
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1349.625
$ws.Range("I81").Value = 1256.7142
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 2513.4284
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -1452.4284
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 1349.625
$ws.Range("I84").Value = 1256.7142
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 12567.142
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -7263.142
$ws.Range("N84").Value = -30608
$ws.Range("H126").Value = 603.5
$ws.Range("I126").Value = 392.77777
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 1178.33331
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = 1291.66669
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 19123.654
$ws.Range("I132").Value = 21061.184
$ws.Range("J132").Value = 3300.5
$ws.Range("K132").Value = 63183.552
$ws.Range("L132").Value = 9901.5
$ws.Range("M132").Value = -60653.552
$ws.Range("N132").Value = -14961.5
$ws.Range("H136").Value = 7694708.5
$ws.Range("I136").Value = 2458.8635
$ws.Range("J136").Value = 23811802
$ws.Range("K136").Value = 7376.5905
$ws.Range("L136").Value = 71435406
$ws.Range("M136").Value = -4826.5905
$ws.Range("N136").Value = -71440506
